# Updated cryptos list on Fri Dec  1 08:58:37 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even when it looks numeric
# (e.g. "61.13", "0.0844") by forcing the cell to a text number-format
# before the assignment, then reverting the style to the default "Normal"
# so no stray style index is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "38.355.63"
$ws.Range("E2").Value = "  +1.66%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.096.77"
$ws.Range("E3").Value = "  +3.43%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "228.44"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.36%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "61.13"
$ws.Range("E7").Value = "  +2.08%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.45%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0844"
$ws.Range("E10").Value = "  +3.00%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.32%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "2.407.54"
$ws.Range("E12").Value = "  +3.52%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +2.76%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "22.32"
$ws.Range("E14").Value = "  +6.29%  "

# Row 15 - Polygon
Set-TextValue $ws.Range("D15") "0.777"
$ws.Range("E15").Value = "  +2.52%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "5.45"
$ws.Range("E16").Value = "  +5.48%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.093.88"
$ws.Range("E17").Value = "  +3.65%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "38.315.15"
$ws.Range("E18").Value = "  +1.76%  "

# Row 19 & 20 - Litecoin and Uniswap swap places
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "6.01"
$ws.Range("E19").Value = "  +1.91%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D20") "70.35"
$ws.Range("E20").Value = "  +1.43%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +1.50%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "225.07"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.06%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.47%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.30"
$ws.Range("E25").Value = "  +2.56%  "

# Row 26 - Monero
Set-TextValue $ws.Range("D26") "169.64"
$ws.Range("E26").Value = "  +1.42%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.40"
$ws.Range("E27").Value = "  +1.30%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.64%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.01"
$ws.Range("E29").Value = "  +1.27%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "1.36"
$ws.Range("E30").Value = "  +8.51%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.120"
$ws.Range("E31").Value = "  -0.04%  "

# Row 32 - WEMIXToken
Set-TextValue $ws.Range("D32") "2.36"
$ws.Range("E32").Value = "  +6.98%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "4.76"
$ws.Range("E33").Value = "  +6.47%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "4.45"
$ws.Range("E34").Value = "  +1.61%  "

# Row 35 - Hedera
Set-TextValue $ws.Range("D35") "0.0604"
$ws.Range("E35").Value = "  +0.10%  "

# Row 36 - THORChain
Set-TextValue $ws.Range("D36") "6.42"
$ws.Range("E36").Value = "  +0.55%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +3.61%  "

# Row 38 - RenderToken
Set-TextValue $ws.Range("D38") "3.49"
$ws.Range("E38").Value = "  +2.96%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  +0.08%  "

# Row 40 - InjectiveProtocol
Set-TextValue $ws.Range("D40") "18.19"
$ws.Range("E40").Value = "  +1.38%  "

# Row 41 - Maker
Set-TextValue $ws.Range("D41") "1.539.88"
$ws.Range("E41").Value = "  +0.54%  "

# Row 42 - Aave
Set-TextValue $ws.Range("D42") "99.84"
$ws.Range("E42").Value = "  +4.58%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  +1.77%  "

# Row 44 - HuobiToken
$ws.Range("E44").Value = "  +1.00%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  +0.17%  "

# Row 46 - FTXToken
$ws.Range("E46").Value = "  +1.77%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +1.03%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +5.68%  "

# Row 49 - ARBITRUM
$ws.Range("E49").Value = "  +3.29%  "

# Row 50 - MXToken
$ws.Range("E50").Value = "  +0.74%  "

# Row 51 - RocketPoolETH
Set-TextValue $ws.Range("D51") "2.293.35"
$ws.Range("E51").Value = "  +3.48%  "
